# Generate Report for Archive
#
# 1. The "Status" value used for both languages on the Overview sheet (and
#    the matching per-language "Status" column on the zh-cn / de-de sheets)
#    moves from "Ready for handoff" to "In Translation".
# 2. The two "Status" columns on the Overview sheet, plus the "Status"
#    column on each language sheet, are narrowed to match the new content.

$wb = $excel.ActiveWorkbook

# --- Update status text -----------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the Status columns ------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
